# Operations Management Assignment.docx - add the SCM reference link
#
# The last paragraph of the document ends with a single space character.
# We insert the slideshare URL (as a real hyperlink) right after that
# existing space, keep a trailing space after the link, and turn the
# inserted URL text into a hyperlink field pointing at the same address
# (so Address == TextToDisplay, matching a plain pasted-link style edit).

$d = $word.ActiveDocument

$url = "https://www.slideshare.net/gadekar1986/supply-chain-management-12816039"

# Last paragraph in the document (currently just a single space).
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)

# Collapsed insertion point right before the paragraph mark, i.e. right
# after the existing trailing space, so that space stays in front of the
# new link.
$insertAt = $d.Range($lastPara.Range.End - 1, $lastPara.Range.End - 1)

$urlStart = $insertAt.Start
$insertAt.InsertAfter($url + " ")
$urlEnd = $urlStart + $url.Length

# Range that covers exactly the freshly inserted URL text (not the
# trailing space we just added).
$urlRange = $d.Range($urlStart, $urlEnd)

# Turn that plain text into a real hyperlink (Address == displayed text).
$d.Hyperlinks.Add($urlRange, $url, $null, $null, $url)
